$d = $word.ActiveDocument

$replacements = @{
  "A imagem apresenta a tela inicial da Google Play Store*" = "A imagem mostra a interface da Google Play Store, na seção ""Para você"". No topo, estão destacados jogos como ""Blood Strike"", ""Roblox"" e ""Tile Club"". Abaixo, há uma lista de sugestões de jogos patrocinados, incluindo ""Paciência"", ""Coin Master"" e ""Bubble Pop! Cannon Shooter"", cada um acompanhado de informações de classificação e tamanho. Na parte inferior da tela, há ícones para navegar entre seções, incluindo ""Jogos"", ""Apps"", ""Livros"" e uma opção de ""Pesquisar"" destacada em vermelho. image_rId8.png"
  "A imagem mostra uma tela de pesquisa em um aplicativo de loja digital*" = "A imagem mostra uma tela de pesquisa na loja de aplicativos, onde o termo ""smart sales force"" está sendo utilizado. Os resultados incluem vários aplicativos, com destaque para ""Smart Força de Vendas"" da Arpa Sistemas, que possui uma classificação de 4,3 estrelas e 14 MB de tamanho, além de mais de mil downloads. Outros aplicativos listados incluem Salesforce, App Sales Force +, e Meta Sales Force, com diferentes classificações e tamanhos. A interface apresenta também um botão de instalação para os aplicativos. image_rId9.png"
  "A imagem apresenta a página de download do aplicativo*" = "A imagem apresenta a interface do aplicativo ""Smart Força de Vendas"", desenvolvido pela Arpa Sistemas. Na parte superior, está o nome do aplicativo junto com a sua classificação de 4,2 estrelas, o número de avaliações (12) e o tamanho do aplicativo (14 MB). Abaixo, há uma chamada para ação para instalar o aplicativo. A imagem exibe também várias capturas de tela do aplicativo, mostrando suas funcionalidades. Há seções como ""Sobre este app"" e ""Segurança dos dados"" também apresentadas na parte inferior. Além disso, são visualizados ícones representando diferentes categorias como jogos, apps, e livros. image_rId10.png"
  "A imagem mostra a tela de instalação do aplicativo*" = "A imagem exibe a tela de instalação do aplicativo ""Smart Força de Vendas"" em um dispositivo móvel. Acima, há um botão para cancelar ou abrir o aplicativo, além de um aviso indicando que ele é verificado pelo Play Protect. Abaixo, são apresentadas sugestões de aplicativos patrocinados, como ""Nomad: Conta em Dólar e Cartão"", ""Livelo: juntar e trocar pontos"" e ""Estoque, Vendas, Pdv, Finanças"", juntamente com mais opções de aplicativos para testar, incluindo ""PictureThis Identificador Planta"" e ""CamScanner"". A parte inferior da tela contém ícones de acesso a jogos, aplicativos, e livros. image_rId11.png"
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    foreach ($pattern in $replacements.Keys) {
        if ($t -like $pattern) {
            $p.Range.Text = $replacements[$pattern]
            break
        }
    }
}
